$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value could be misread as a number by Excel;
# force text format first so the value round-trips as a string.
$textCells = @("D5", "D6", "D7", "D9", "D11", "D12", "D15", "D16", "D20", "D21", "D22", "D26", "D28", "D29", "D31", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '44.140.47'
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').Value = '2.355.63'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '241.39'
$ws.Range('E5').Value = '  +3.33%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').Value = '0.671'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('D7').Value = '74.28'
$ws.Range('E7').Value = '  +6.37%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '0.569'
$ws.Range('E9').Value = '  +23.20%  '
$ws.Range('E10').Value = '  +5.03%  '
$ws.Range('D11').Value = '31.85'
$ws.Range('E11').Value = '  +21.15%  '
$ws.Range('D12').Value = '7.38'
$ws.Range('E12').Value = '  +18.25%  '
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').Value = '2.704.41'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').Value = '16.82'
$ws.Range('E15').Value = '  +7.04%  '
$ws.Range('D16').Value = '0.908'
$ws.Range('E16').Value = '  +6.38%  '
$ws.Range('D17').Value = '2.360.32'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '44.257.75'
$ws.Range('E19').Value = '  +3.97%  '
$ws.Range('D20').Value = '6.70'
$ws.Range('E20').Value = '  +5.65%  '
$ws.Range('D21').Value = '78.13'
$ws.Range('E21').Value = '  +5.42%  '
$ws.Range('D22').Value = '255.70'
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -4.57%  '
$ws.Range('E25').Value = '  +3.52%  '
$ws.Range('D26').Value = '10.71'
$ws.Range('E26').Value = '  +7.16%  '
$ws.Range('E27').Value = '  +3.65%  '
$ws.Range('D28').Value = '22.57'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('D29').Value = '174.80'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('D31').Value = '0.130'
$ws.Range('E31').Value = '  +3.24%  '
$ws.Range('E32').Value = '  +4.89%  '
$ws.Range('D33').Value = '5.39'
$ws.Range('E33').Value = '  +8.15%  '
$ws.Range('D34').Value = '0.0757'
$ws.Range('E34').Value = '  +9.51%  '
$ws.Range('D35').Value = '5.33'
$ws.Range('E35').Value = '  +4.74%  '
$ws.Range('E36').Value = '  +6.25%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '2.45'
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').Value = '6.55'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('E39').Value = '  +7.24%  '
$ws.Range('D40').Value = '19.28'
$ws.Range('E40').Value = '  +3.63%  '
$ws.Range('D41').Value = '8.98'
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '0.191'
$ws.Range('E43').Value = '  +14.20%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '1.26'
$ws.Range('E44').Value = '  +3.00%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = '0.0994'
$ws.Range('E45').Value = '  +5.03%  '
$ws.Range('D46').Value = '2.48'
$ws.Range('E46').Value = '  +10.55%  '
$ws.Range('D47').Value = '100.60'
$ws.Range('E47').Value = '  +1.17%  '
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('D49').Value = '4.47'
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').Value = '0.000209'
$ws.Range('E50').Value = '  +4.22%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.451.45'
$ws.Range('E51').Value = '  +0.04%  '

# Remove the temporary text formatting so cell styling matches the original
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
